$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in D/E are kept as text while we set them,
# matching the original inline-string cell content (avoids float coercion).
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.375.84"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "3.400.69"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "581.36"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").Value = "179.25"
$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("E9").Value = "  +8.21%  "

$ws.Range("D10").Value = "0.587"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").Value = "48.44"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "0.0000284"
$ws.Range("E12").Value = "  +3.68%  "

$ws.Range("D13").Value = "681.52"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.947.37"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "69.462.32"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").Value = "3.399.51"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "17.70"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("D20").Value = "11.29"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").Value = "0.910"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("E22").Value = "  -1.59%  "

$ws.Range("D23").Value = "17.09"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").Value = "101.18"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").Value = "3.90"
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("E27").Value = "  +2.46%  "

$ws.Range("D28").Value = "33.55"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("D29").Value = "8.77"
$ws.Range("E29").Value = "  +2.78%  "

$ws.Range("E30").Value = "  -0.75%  "

$ws.Range("D31").Value = "3.79"
$ws.Range("E31").Value = "  +13.46%  "

$ws.Range("D32").Value = "556.48"
$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").Value = "58.03"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "3.609.03"
$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("D39").Value = "35.31"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  +11.42%  "

$ws.Range("E41").Value = "  +4.32%  "

$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +3.62%  "

$ws.Range("E43").Value = "  +3.89%  "

$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").Value = "1.40"
$ws.Range("E47").Value = "  +4.17%  "

$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("D49").Value = "131.09"

$ws.Range("D50").Value = "2.63"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("D51").Value = "7.42"
$ws.Range("E51").Value = "  -0.38%  "

# Restore default styling so only cell contents differ from the original file.
$priceRange.Style = "Normal"
